$d = $word.ActiveDocument

# --- Change 1: merge the JFrame/JPanel runs in the graphics-stub comment paragraph ---
$d.Content.Find.Execute(" or JPanel, give it appropriate size,", $true, $false, $false, $false, $false, $true, 1, $false, " or JPanel, give it appropriate size,", 2) | Out-Null

# --- Change 2: flesh out the addQueens() pseudocode before the trailing bookmark paragraph ---
# Locate the blank, left-indented (720 twips) paragraph that sits right before the blank,
# non-indented separator paragraph near the end of the EightQueens class body. That is the
# paragraph which (invisibly) carries the _GoBack bookmark.
$cr = [char]13
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Format.LeftIndent -eq 36 -and $p.Range.Text -eq $cr) {
        $next = $p.Next()
        if ($next -ne $null -and $next.Format.LeftIndent -eq 0 -and $next.Range.Text -eq $cr) {
            $target = $p
            break
        }
    }
}

# Drop the paragraph's left indent (720 -> none) to match the new formatting.
$target.Format.LeftIndent = 0

$insertXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>public</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>boolean</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>addQueens</w:t>
      </w:r>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>ArrayList</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">&lt;Queen&gt; </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>alreadyPlaced</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>) {</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>if</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> (</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>alreadyPlaced.size</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">() &gt;= </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>ROWS</w:t>
      </w:r>
      <w:r>
        <w:t>) {</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>// base case - if there are 8 queens on the board</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>return</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>true</w:t>
      </w:r>
      <w:r>
        <w:t>;</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t xml:space="preserve">} </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>else</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>if</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> (</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>alreadyPlaced.isEmpty</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>()) {</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>// if the array is empty, place a queen in the start position, which changes</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>// based on the solution number, and then call addQueens again</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>Queen</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>q</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> = </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>new</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>Queen(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>startRow</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>startCol</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>);</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>alreadyPlaced.add</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>(q);</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>return</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> addQueens(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>alreadyPlaced</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>);</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t xml:space="preserve">} </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>else</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> {</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>int</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>curSize</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> = </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>alreadyPlaced.size</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>();</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>//a quick instance field to check if a Queen has been added</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> later</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:ind w:left="720" w:firstLine="720"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>for</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> (</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>int</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>i</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> = 0; </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>i</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> &lt; </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>COLUMNS;</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>i</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>++) {</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:ind w:left="2160"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>Queen</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>q</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> = </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>new</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Queen(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>alreadyPlaced.get</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>alreadyPlaced.size</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">() </w:t>
      </w:r>
      <w:r>
        <w:t>-</w:t>
      </w:r>
      <w:r>
        <w:t>1).</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>getRow</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>()</w:t>
      </w:r>
      <w:r>
        <w:t>+1</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>i</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>);</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:ind w:left="2160"/>
      </w:pPr>
      <w:r>
        <w:t>//starting on the next row, try placing a Queen in each spot. If one fits, then break the loop</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>if</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> (</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>isLegal</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve">q, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>alreadyPlaced</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>)) {</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>alreadyPlaced.add</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>(q);</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>break</w:t>
      </w:r>
      <w:r>
        <w:t>;</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>if</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> (</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>curSize</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> &lt; </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>alreadyPlaced.size</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>()) {</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>return</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> addQueens(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>alreadyPlaced</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>);</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
        <w:t>//if a Queen has been added, t</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:t>hat means there is potential for more, so call the method again</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t xml:space="preserve">} </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>else</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>return</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>false</w:t>
      </w:r>
      <w:r>
        <w:t>;</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:ind w:firstLine="720"/>
      </w:pPr>
      <w:r>
        <w:t>//if no Queen has been added, terminate the method and return false</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
        <w:t>}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>}</w:t>
      </w:r>
    </w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertRange = $d.Range($target.Range.Start, $target.Range.Start)
$insertRange.InsertXML($insertXml)
